$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Winglets" question column (B) -------------------------------
# B1 header, matching the centered style already used by A1/A2
$ws.Range("B1").Value = "Winglets"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# B2 answer value
$ws.Range("B2").Value = "Yes"

# Give column B an explicit width
$ws.Columns("B").ColumnWidth = 14

# Data validation drop-down list on B2
$val = $ws.Range("B2").Validation
$val.Add(3, 1, 1, """Yes,No,Don't Know""")
$val.InputTitle = "Winglets"
$val.InputMessage = "Are Winglets installed: Yes / No"
$val.ShowInput = $true
$val.ShowError = $true

# Restore the selection Excel left the cursor on
[void]$ws.Range("F5").Select()
